$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-5 with the new TPM-derived values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icosl"
$ws.Range("C2").Value = "Cd28"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9314806666666667
$ws.Range("H2").Value = 2.794442
$ws.Range("I2").Value = 0.06288211749152639
$ws.Range("J2").Value = 0.06288211749152638
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.961972333333333
$ws.Range("N2").Value = 14.885917
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 4.621981297034889
$ws.Range("R2").Value = 41.597831673314
$ws.Range("S2").Value = 0.06288211749152639
$ws.Range("T2").Value = 0.06288211749152638

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Icosl"
$ws.Range("C3").Value = "Cd28"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.399281333333333
$ws.Range("H3").Value = 7.197844
$ws.Range("I3").Value = 0.1619699646990985
$ws.Range("J3").Value = 0.1619699646990985
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.961972333333333
$ws.Range("N3").Value = 14.885917
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 11.90516759588311
$ws.Range("R3").Value = 107.146508362948
$ws.Range("S3").Value = 0.1619699646990985
$ws.Range("T3").Value = 0.1619699646990985

$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Icosl"
$ws.Range("C4").Value = "Cd28"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.416663
$ws.Range("H4").Value = 1.249989
$ws.Range("I4").Value = 0.02812796084553394
$ws.Range("J4").Value = 0.02812796084553394
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.961972333333333
$ws.Range("N4").Value = 14.885917
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 2.067470278323666
$ws.Range("R4").Value = 18.607232504913
$ws.Range("S4").Value = 0.02812796084553394
$ws.Range("T4").Value = 0.02812796084553394

$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Icosl"
$ws.Range("C5").Value = "Cd28"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.0657
$ws.Range("H5").Value = 33.1971
$ws.Range("I5").Value = 0.7470199569638412
$ws.Range("J5").Value = 0.7470199569638412
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.961972333333333
$ws.Range("N5").Value = 14.885917
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 54.90769724896666
$ws.Range("R5").Value = 494.1692752407
$ws.Range("S5").Value = 0.7470199569638412
$ws.Range("T5").Value = 0.7470199569638412

# Remove the now-obsolete rows 6-9 (MuSCs/Resolving-Mac combos collapsed away)
$ws.Range("A6:T9").EntireRow.Delete()
